$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fbn1"
$ws.Range("C2").Value = "Itgav"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 7.848425
$ws.Range("H2").Value = 23.545275
$ws.Range("I2").Value = 0.02436729568045431
$ws.Range("J2").Value = 0.02436729568045431
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 16.535604
$ws.Range("N2").Value = 49.606812
$ws.Range("O2").Value = 0.2120453146491552
$ws.Range("P2").Value = 0.2120453146491552
$ws.Range("Q2").Value = 129.7784478237
$ws.Range("R2").Value = 1168.0060304133
$ws.Range("S2").Value = 0.005166970879710934
$ws.Range("T2").Value = 0.005166970879710935

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fbn1"
$ws.Range("C3").Value = "Itgav"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 7.848425
$ws.Range("H3").Value = 23.545275
$ws.Range("I3").Value = 0.02436729568045431
$ws.Range("J3").Value = 0.02436729568045431
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 40.62063066666667
$ws.Range("N3").Value = 121.861892
$ws.Range("O3").Value = 0.5209011059384622
$ws.Range("P3").Value = 0.5209011059384622
$ws.Range("Q3").Value = 318.8079732400333
$ws.Range("R3").Value = 2869.2717591603
$ws.Range("S3").Value = 0.01269295126867816
$ws.Range("T3").Value = 0.01269295126867816

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fbn1"
$ws.Range("C4").Value = "Itgav"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 7.848425
$ws.Range("H4").Value = 23.545275
$ws.Range("I4").Value = 0.02436729568045431
$ws.Range("J4").Value = 0.02436729568045431
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 20.825229
$ws.Range("N4").Value = 62.475687
$ws.Range("O4").Value = 0.2670535794123827
$ws.Range("P4").Value = 0.2670535794123827
$ws.Range("Q4").Value = 163.445247914325
$ws.Range("R4").Value = 1471.007231228925
$ws.Range("S4").Value = 0.006507373532065213
$ws.Range("T4").Value = 0.006507373532065213

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fbn1"
$ws.Range("C5").Value = "Itgav"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 291.329961
$ws.Range("H5").Value = 873.989883
$ws.Range("I5").Value = 0.9045029162236017
$ws.Range("J5").Value = 0.9045029162236017
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 16.535604
$ws.Range("N5").Value = 49.606812
$ws.Range("O5").Value = 0.2120453146491552
$ws.Range("P5").Value = 0.2120453146491552
$ws.Range("Q5").Value = 4817.316868431443
$ws.Range("R5").Value = 43355.85181588299
$ws.Range("S5").Value = 0.1917956054717121
$ws.Range("T5").Value = 0.1917956054717121

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fbn1"
$ws.Range("C6").Value = "Itgav"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 291.329961
$ws.Range("H6").Value = 873.989883
$ws.Range("I6").Value = 0.9045029162236017
$ws.Range("J6").Value = 0.9045029162236017
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 40.62063066666667
$ws.Range("N6").Value = 121.861892
$ws.Range("O6").Value = 0.5209011059384622
$ws.Range("P6").Value = 0.5209011059384622
$ws.Range("Q6").Value = 11834.0067479154
$ws.Range("R6").Value = 106506.0607312386
$ws.Range("S6").Value = 0.4711565693854383
$ws.Range("T6").Value = 0.4711565693854383

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fbn1"
$ws.Range("C7").Value = "Itgav"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 291.329961
$ws.Range("H7").Value = 873.989883
$ws.Range("I7").Value = 0.9045029162236017
$ws.Range("J7").Value = 0.9045029162236017
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 20.825229
$ws.Range("N7").Value = 62.475687
$ws.Range("O7").Value = 0.2670535794123827
$ws.Range("P7").Value = 0.2670535794123827
$ws.Range("Q7").Value = 6067.013152386068
$ws.Range("R7").Value = 54603.11837147462
$ws.Range("S7").Value = 0.2415507413664513
$ws.Range("T7").Value = 0.2415507413664513

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fbn1"
$ws.Range("C8").Value = "Itgav"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 22.91008466666667
$ws.Range("H8").Value = 68.730254
$ws.Range("I8").Value = 0.07112978809594397
$ws.Range("J8").Value = 0.07112978809594397
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 16.535604
$ws.Range("N8").Value = 49.606812
$ws.Range("O8").Value = 0.2120453146491552
$ws.Range("P8").Value = 0.2120453146491552
$ws.Range("Q8").Value = 378.832087654472
$ws.Range("R8").Value = 3409.488788890248
$ws.Range("S8").Value = 0.01508273829773218
$ws.Range("T8").Value = 0.01508273829773218

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fbn1"
$ws.Range("C9").Value = "Itgav"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 22.91008466666667
$ws.Range("H9").Value = 68.730254
$ws.Range("I9").Value = 0.07112978809594397
$ws.Range("J9").Value = 0.07112978809594397
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 40.62063066666667
$ws.Range("N9").Value = 121.861892
$ws.Range("O9").Value = 0.5209011059384622
$ws.Range("P9").Value = 0.5209011059384622
$ws.Range("Q9").Value = 930.6220877867298
$ws.Range("R9").Value = 8375.598790080569
$ws.Range("S9").Value = 0.03705158528434568
$ws.Range("T9").Value = 0.03705158528434568

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fbn1"
$ws.Range("C10").Value = "Itgav"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 22.91008466666667
$ws.Range("H10").Value = 68.730254
$ws.Range("I10").Value = 0.07112978809594397
$ws.Range("J10").Value = 0.07112978809594397
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 20.825229
$ws.Range("N10").Value = 62.475687
$ws.Range("O10").Value = 0.2670535794123827
$ws.Range("P10").Value = 0.2670535794123827
$ws.Range("Q10").Value = 477.107759592722
$ws.Range("R10").Value = 4293.969836334498
$ws.Range("S10").Value = 0.01899546451386613
$ws.Range("T10").Value = 0.01899546451386613
